$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing the existing rows 7-44 down to 8-45.
$ws.Rows(7).Insert()

# Populate the new row 7 with a new weekly price observation. Most columns
# (market/region/category/unit/quality/origin/classification) are identical
# to the rest of the sheet's rows; only the date, volume and weighted price
# are new for this entry.
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44503
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 100112001
$ws.Cells.Item(7, 7).Value = "Berenjena"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 130
$ws.Cells.Item(7, 11).Value = 8000
$ws.Cells.Item(7, 12).Value = 9000
$ws.Cells.Item(7, 13).Value = 8538
$ws.Cells.Item(7, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 142
$ws.Cells.Item(7, 17).Value = 60
$ws.Cells.Item(7, 18).Value = "Hortaliza"
